# This script applies the "corrections create results and IRR calculation" edit:
# adds a new results row (index 3 / year 2020) to both the "Overview" and
# "Capacity" worksheets.

$wb = $excel.ActiveWorkbook

# ----- Sheet "Overview" -----
$ws1 = $wb.Worksheets.Item("Overview")

# Copy the style of A4 (bold/bordered header-like style) down into A5, then
# set its value to the new row index.
$ws1.Range("A4").Copy($ws1.Range("A5"))
$ws1.Range("A5").Value = 3

$ws1.Cells.Item(5, 2).Value = 2020
$ws1.Cells.Item(5, 3).Value = 544999999.1
$ws1.Cells.Item(5, 4).Value = 20811313152.32112
$ws1.Cells.Item(5, 5).Value = 38.1858957553916
$ws1.Cells.Item(5, 6).Value = 11
$ws1.Cells.Item(5, 7).Value = 222307.2994627971
$ws1.Cells.Item(5, 8).Value = 0
$ws1.Cells.Item(5, 9).Value = 2.560297453063297
$ws1.Cells.Item(5, 10).Value = 97504.61057238668
$ws1.Cells.Item(5, 11).Value = 17769.20341851865
$ws1.Cells.Item(5, 12).Value = 0.1822396224568984
$ws1.Cells.Item(5, 13).Value = 1
$ws1.Cells.Item(5, 14).Value = 0
$ws1.Cells.Item(5, 15).Value = 711988264.0362152
$ws1.Cells.Item(5, 16).Value = 0

# ----- Sheet "Capacity" -----
$ws2 = $wb.Worksheets.Item("Capacity")

# Copy the style of A4 into A5, then set its value to the new row index.
$ws2.Range("A4").Copy($ws2.Range("A5"))
$ws2.Range("A5").Value = 3

$ws2.Cells.Item(5, 2).Value = 2020
$ws2.Cells.Item(5, 3).Value = 4644.4034
$ws2.Cells.Item(5, 4).Value = 25208582.8382924
$ws2.Cells.Item(5, 5).Value = 954932128.661841
$ws2.Cells.Item(5, 6).Value = 37.88123016623044
$ws2.Cells.Item(5, 7).Value = 24845.77
$ws2.Cells.Item(5, 8).Value = 59003.61621933627
$ws2.Cells.Item(5, 9).Value = 3332129.912351787
$ws2.Cells.Item(5, 10).Value = 56.47331682121212
$ws2.Cells.Item(5, 11).Value = 31358.329
$ws2.Cells.Item(5, 12).Value = 217108263.1152519
$ws2.Cells.Item(5, 13).Value = 8460882562.993628
$ws2.Cells.Item(5, 14).Value = 38.97079936797324
$ws2.Cells.Item(5, 15).Value = 8194.3025
$ws2.Cells.Item(5, 16).Value = 1007081.087896536
$ws2.Cells.Item(5, 17).Value = 53413355.2968122
$ws2.Cells.Item(5, 18).Value = 53.03779004367492
$ws2.Cells.Item(5, 19).Value = 8858.749999999998
$ws2.Cells.Item(5, 20).Value = 18624635.99999999
$ws2.Cells.Item(5, 21).Value = 705091697.2187846
$ws2.Cells.Item(5, 22).Value = 37.85801221665675
$ws2.Cells.Item(5, 23).Value = 8599
$ws2.Cells.Item(5, 24).Value = 74259244.44886312
$ws2.Cells.Item(5, 25).Value = 2846665877.166388
$ws2.Cells.Item(5, 26).Value = 38.33416160228613
$ws2.Cells.Item(5, 27).Value = 47547.50848700004
$ws2.Cells.Item(5, 28).Value = 80823362.9723005
$ws2.Cells.Item(5, 29).Value = 2945418492.530766
$ws2.Cells.Item(5, 30).Value = 36.44266192610929
$ws2.Cells.Item(5, 31).Value = 10271.8
$ws2.Cells.Item(5, 32).Value = 40873004.28171189
$ws2.Cells.Item(5, 33).Value = 1517259559.244995
$ws2.Cells.Item(5, 34).Value = 37.12131236518557
$ws2.Cells.Item(5, 35).Value = 53555.51607579708
$ws2.Cells.Item(5, 36).Value = 50877740.2720072
$ws2.Cells.Item(5, 37).Value = 1831210598.874818
$ws2.Cells.Item(5, 38).Value = 35.9923728743579
$ws2.Cells.Item(5, 39).Value = 20779.02
$ws2.Cells.Item(5, 40).Value = 36159080.46745713
$ws2.Cells.Item(5, 41).Value = 1493106750.420732
$ws2.Cells.Item(5, 42).Value = 41.29271903815462
$ws2.Cells.Item(5, 43).Value = 3652.9
$ws2.Cells.Item(5, 44).Value = 0
$ws2.Cells.Item(5, 45).Value = 0
$ws2.Cells.Item(5, 46).Value = 0
$ws2.Cells.Item(5, 47).Value = 0
$ws2.Cells.Item(5, 48).Value = 0
$ws2.Cells.Item(5, 49).Value = 0
$ws2.Cells.Item(5, 50).Value = 0
$ws2.Cells.Item(5, 51).Value = 0
$ws2.Cells.Item(5, 52).Value = 0
$ws2.Cells.Item(5, 53).Value = 0
$ws2.Cells.Item(5, 54).Value = 0
